# Add summary tag with label in group feature
# The <<sum>> tag in G6 becomes &="Total: "<<sum>>
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = '&="Total: "<<sum>>'
